# "last commit of the day" - fill in the remaining BOM rows (21-27) of the
# Stückliste with part data, mark their supplier as "Mouser", set the
# corresponding Lieferant/Mouser on row 20 as well, resize the rows that now
# wrap onto multiple lines, and move the active selection to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 only needed the "Lieferant" (Mouser) filled in - Bezeichnung,
# Details and Best.-Nr. were already populated.
$ws.Range("E20").Value = "Mouser"

# Row 21 - 330uF Tantalum capacitor
$ws.Range("B21").Value = "330uF Kondensator"
$ws.Range("C21").Value = "Tantalkondensatoren – Polymer 16V 330uF 7343 ESR=15mOhms 20% "
$ws.Range("E21").Value = "Mouser"
$ws.Range("F21").Value = "80-T521X337M16ATE015 "
$ws.Rows("21").RowHeight = 49.5

# Row 22 - 187uH inductor
$ws.Range("B22").Value = "187 uH Spule"
$ws.Range("C22").Value = "Festinduktivitäten 187uH 8A 20% DCR=37.6mOhms "
$ws.Range("E22").Value = "Mouser"
$ws.Range("F22").Value = "80-SHBC14-1R2A0187V "
$ws.Rows("22").RowHeight = 33

# Row 23 - 4.7mF electrolytic capacitor
$ws.Range("B23").Value = "4.7mF Kondensator"
$ws.Range("C23").Value = "Aluminium-Elektrolyt-Kondensatoren - Radial bedrahtet 4700uF 10V 20% LYTICS/IC "
$ws.Range("E23").Value = "Mouser"
$ws.Range("F23").Value = "598-478CKE010M "
$ws.Rows("23").RowHeight = 66

# Row 24 - 68k resistor (no Details/Lieferant description column value)
$ws.Range("B24").Value = "68k Widerstand"
$ws.Range("E24").Value = "Mouser"
$ws.Range("F24").Value = "603-RC0805FR-0768KL"

# Row 25 - 6k8 resistor
$ws.Range("B25").Value = "6k8 Widerstand"
$ws.Range("E25").Value = "Mouser"
$ws.Range("F25").Value = "603-RC0805FR-076K8L"

# Row 26 - 2-position terminal block
$ws.Range("B26").Value = "Feste Anschlussblöcke 2P"
$ws.Range("C26").Value = "Feste Anschlussblöcke 2 POS. PC/MNT. TERM. "
$ws.Range("E26").Value = "Mouser"
$ws.Range("F26").Value = "571-2828412 "
$ws.Rows("26").RowHeight = 33

# Row 27 - 3-position terminal block
$ws.Range("B27").Value = "Feste Anschlussblöcke 3P"
$ws.Range("C27").Value = "Feste Anschlussblöcke 3P PCB TERM BLK"
$ws.Range("E27").Value = "Mouser"
$ws.Range("F27").Value = "571-2828413 "
$ws.Rows("27").RowHeight = 33

# Move the selection to where editing left off.
$ws.Range("F11").Select()
